# Work-Estimate.xlsx sample-data cleanup
# The "Example" template sheet had placeholder/sample values filled in for
# demonstration purposes. This commit clears those sample values back out,
# leaving the form blank (ready to be filled in programmatically by the
# WorkEstimateExcelDTO xlsx writer being developed), while keeping the
# genuine label/template text (e.g. "Gestia Armatora", "Repair handling")
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: sample container prefix / serial / check-digit
$ws.Range("H7").Value = $null
$ws.Range("J7").Value = $null
$ws.Range("M7").Value = $null

# Row 9: sample in/out date
$ws.Range("H9").Value = $null
# Row 9: sample container type
$ws.Range("K9").Value = $null

# Row 10: sample terminal name
$ws.Range("B10").Value = $null
# Row 10: sample "date of estimate"
$ws.Range("D10").Value = $null

# Sample numeric entries in the repair-lines table that fed the totals
$ws.Range("M18").Value = $null
$ws.Range("M21").Value = $null

# Remove the LOCONI logo picture that used to sit over B4:E6
if ($ws.Shapes.Count -gt 0) {
    for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
        $ws.Shapes.Item($i).Delete()
    }
}
